# FieldRenameMap.xlsx edit:
#  - Sheet1!C2:C7 ("ENABLED" column for the Variable-related rows) flips
#    from "Y" to "N".
#  - The saved selection/scroll position moves to E8 (no frozen
#    topLeftCell override).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = "N"

$ws.Range("E8").Select()
